$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.301.23"
$ws.Range("E2").Value = "  -0.04%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.932.54"
$ws.Range("E3").Value = "  -0.02%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9981"
$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7848"
$ws.Range("E5").Value = "  +9.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "245.58"
$ws.Range("E6").Value = "  -2.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9964"
$ws.Range("E7").Value = "  -0.39%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3215"
$ws.Range("E8").Value = "  -2.33%  "

$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.93"
$ws.Range("E9").Value = "  +1.52%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07083"
$ws.Range("E10").Value = "  -2.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7817"
$ws.Range("E11").Value = "  -2.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07993"
$ws.Range("E12").Value = "  -1.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.930.86"
$ws.Range("E13").Value = "  -0.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.368"
$ws.Range("E14").Value = "  -1.72%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.84"
$ws.Range("E15").Value = "  +0.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.54"
$ws.Range("E16").Value = "  -3.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.290.67"
$ws.Range("E17").Value = "  -0.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "255.61"
$ws.Range("E18").Value = "  +1.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007999"
$ws.Range("E19").Value = "  -2.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.811"
$ws.Range("E20").Value = "  +0.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.188.94"
$ws.Range("E21").Value = "  +0.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9973"
$ws.Range("E22").Value = "  -0.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9981"
$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.768"
$ws.Range("E24").Value = "  -2.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.583"
$ws.Range("E25").Value = "  -1.59%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.50"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1368"
$ws.Range("E27").Value = "  +5.72%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.09"
$ws.Range("E28").Value = "  -1.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.277"
$ws.Range("E29").Value = "  -2.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.368"
$ws.Range("E30").Value = "  +1.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.514"
$ws.Range("E31").Value = "  -1.90%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.422"
$ws.Range("E32").Value = "  +0.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.137"
$ws.Range("E33").Value = "  -0.91%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05177"
$ws.Range("E34").Value = "  -0.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.285"
$ws.Range("E35").Value = "  +2.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7516"
$ws.Range("E36").Value = "  +0.83%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.765"
$ws.Range("E37").Value = "  -0.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01964"
$ws.Range("E38").Value = "  -0.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.799"
$ws.Range("E39").Value = "  -0.28%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "78.86"
$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.412"
$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4517"
$ws.Range("E42").Value = "  -0.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.977"
$ws.Range("E43").Value = "  -1.94%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9971"
$ws.Range("E44").Value = "  -0.31%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8334"
$ws.Range("E45").Value = "  -1.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.82"
$ws.Range("E46").Value = "  -0.57%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.788"
$ws.Range("E47").Value = "  +0.33%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.502"
$ws.Range("E48").Value = "  +0.89%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "982.36"
$ws.Range("E49").Value = "  +11.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.21"
$ws.Range("E50").Value = "  +1.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4167"
$ws.Range("E51").Value = "  +0.02%  "
